$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# "update database and change read_price algorithm"
#
# The quarterly/annual period columns roll forward by one: the oldest
# period (column D) is dropped, every remaining period shifts one
# column to the left (D<-E, E<-F, ... L<-M), and the newest period's
# figures land in column M.
#
# xlPasteValues (-4163) is used for every shift so that:
#   - date-like text such as "1401-04-26" / "1401-10-28" is copied
#     verbatim instead of being re-interpreted as a date serial by the
#     normal .Value input-parser, and
#   - the destination cells keep their original style index (no new
#     numFmt/style records get created).
# -----------------------------------------------------------------

$xlPasteValues = -4163

$dataRows = @(8, 9, 11, 12, 13, 14, 16, 17, 19, 20, 21, 22, 24, 26)

# Row 9's old J (col 10) value "1401-10-28 (6)" becomes "1402-02-10 (7)"
# before the shift (it lands in column I once everything moves left).
$ws.Cells.Item(9, 10).Value = "1402-02-10 (7)"

foreach ($r in $dataRows) {
    $ws.Range($ws.Cells.Item($r, 5), $ws.Cells.Item($r, 13)).Copy() | Out-Null
    $ws.Range($ws.Cells.Item($r, 4), $ws.Cells.Item($r, 12)).PasteSpecial($xlPasteValues) | Out-Null
}
$excel.CutCopyMode = 0

# Newest period (column M) values.
$ws.Cells.Item(8, 13).Value  = "12 ماهه منتهی به 1401/12"

# "1402-02-10" looks like a date, so route it through a scratch-cell
# formula + paste-values instead of a direct .Value assignment (which
# Excel's input parser would silently turn into a date serial number).
$scratch = $ws.Cells.Item(200, 200)
$scratch.Formula = '="1402-02-10"'
$scratch.Copy() | Out-Null
$ws.Cells.Item(9, 13).PasteSpecial($xlPasteValues) | Out-Null
$scratch.Clear() | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(11, 13).Value = 224892
$ws.Cells.Item(12, 13).Value = -70152
$ws.Cells.Item(13, 13).Value = 154740
$ws.Cells.Item(14, 13).Value = -14043
$ws.Cells.Item(16, 13).Value = 5079
$ws.Cells.Item(17, 13).Value = 145776
$ws.Cells.Item(19, 13).Value = 11437
$ws.Cells.Item(20, 13).Value = 157213
$ws.Cells.Item(21, 13).Value = -17026
$ws.Cells.Item(22, 13).Value = 140187
$ws.Cells.Item(24, 13).Value = 140187
$ws.Cells.Item(26, 13).Value = 5072

# -----------------------------------------------------------------
# Column widths also roll forward with the data: the "wide" columns
# (29 chars, used for the last period of every fiscal year) move from
# F/J to E/I, and the new rightmost column M becomes wide as well.
# -----------------------------------------------------------------
$narrowWidth = $ws.Columns.Item(4).ColumnWidth   # existing 28-char column
$wideWidth   = $ws.Columns.Item(6).ColumnWidth   # existing 29-char column

$ws.Columns.Item(4).ColumnWidth  = $narrowWidth
$ws.Columns.Item(5).ColumnWidth  = $wideWidth
$ws.Columns.Item(6).ColumnWidth  = $narrowWidth
$ws.Columns.Item(7).ColumnWidth  = $narrowWidth
$ws.Columns.Item(8).ColumnWidth  = $narrowWidth
$ws.Columns.Item(9).ColumnWidth  = $wideWidth
$ws.Columns.Item(10).ColumnWidth = $narrowWidth
$ws.Columns.Item(11).ColumnWidth = $narrowWidth
$ws.Columns.Item(12).ColumnWidth = $narrowWidth
$ws.Columns.Item(13).ColumnWidth = $wideWidth
